$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A1').Value = 'admin'
$ws.Range('C1').Value = 'subAdmin'
$ws.Range('E1').Value = 'customer'
$ws.Range('G1').Value = 'category'
$ws.Range('I1').Value = 'level'
$ws.Range('A2').Value = 'id'
$ws.Range('C2').Value = 'id'
$ws.Range('E2').Value = 'id'
$ws.Range('G2').Value = 'id'
$ws.Range('I2').Value = 'id'
$ws.Range('A3').Value = 'name'
$ws.Range('C3').Value = 'name'
$ws.Range('E3').Value = 'firstName'
$ws.Range('G3').Value = 'name'
$ws.Range('I3').Value = 'name'
$ws.Range('A4').Value = 'password'
$ws.Range('C4').Value = 'password'
$ws.Range('E4').Value = 'lastName'
$ws.Range('E5').Value = 'userName'
$ws.Range('E6').Value = 'phone'
$ws.Range('E7').Value = 'email'
$ws.Range('E8').Value = 'address'
$ws.Range('C10').Value = 'product'
$ws.Range('E10').Value = 'cart'
$ws.Range('G10').Value = 'comment'
$ws.Range('I10').Value = 'rating'
$ws.Range('K10').Value = 'payment'
$ws.Range('M10').Value = 'shipment'
$ws.Range('O10').Value = 'order'
$ws.Range('C11').Value = 'id'
$ws.Range('E11').Value = 'id'
$ws.Range('G11').Value = 'id'
$ws.Range('I11').Value = 'id'
$ws.Range('K11').Value = 'id'
$ws.Range('M11').Value = 'id'
$ws.Range('O11').Value = 'id'
$ws.Range('C12').Value = 'name'
$ws.Range('E12').Value = 'totalCost'
$ws.Range('G12').Value = 'heading'
$ws.Range('I12').Value = 'rate'
$ws.Range('K12').Value = 'amount'
$ws.Range('M12').Value = 'DVtime'
$ws.Range('O12').Value = 'status'
$ws.Range('C13').Value = 'describe'
$ws.Range('E13').Value = 'quantity'
$ws.Range('G13').Value = 'content'
$ws.Range('K13').Value = 'type'
$ws.Range('M13').Value = 'SDTime'
$ws.Range('O13').Value = 'cost'
$ws.Range('C14').Value = 'price'
$ws.Range('E14').Value = 'ten'
$ws.Range('G14').Value = 'time'
$ws.Range('O14').Value = 'size'
$ws.Range('C15').Value = 'img'
$ws.Range('E15').Value = 'sdt '
$ws.Range('O15').Value = 'address'
$ws.Range('C16').Value = 'quantity'
$ws.Range('E16').Value = 'ghi chu'
$ws.Range('O16').Value = 'paymentType'
$ws.Range('C17').Value = 'origin'
$ws.Range('C18').Value = 'status'

$headerCells = @('A1','C1','E1','G1','I1','C10','E10','G10','I10','K10','M10','O10')
foreach ($addr in $headerCells) {
    $ws.Range($addr).Interior.Color = 65535
}

[void]$ws.Range("O23").Select()
